$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 18.96118728265645
$ws.Cells.Item(2, 3).Value = 8.965750565538047
$ws.Cells.Item(2, 4).Value = 8.757245386148615
$ws.Cells.Item(2, 6).Value = 35.16409869795943
$ws.Cells.Item(2, 7).Value = 3.666717079945277
$ws.Cells.Item(2, 10).Value = 10.38146579863548
$ws.Cells.Item(2, 12).Value = 11.7736932647385
$ws.Cells.Item(2, 14).Value = 18.40874426480805
$ws.Cells.Item(2, 15).Value = 26.68862419243844
$ws.Cells.Item(3, 2).Value = 18.55484199323308
$ws.Cells.Item(3, 3).Value = 8.767467454062009
$ws.Cells.Item(3, 4).Value = 8.75863141627808
$ws.Cells.Item(3, 6).Value = 35.21372646509152
$ws.Cells.Item(3, 7).Value = 3.668978235981334
$ws.Cells.Item(3, 10).Value = 10.40727553765922
$ws.Cells.Item(3, 12).Value = 11.75991285292342
$ws.Cells.Item(3, 14).Value = 18.46455057013625
$ws.Cells.Item(3, 15).Value = 26.74725699137069
$ws.Cells.Item(4, 2).Value = 18.30380673043957
$ws.Cells.Item(4, 3).Value = 8.64244192612008
$ws.Cells.Item(4, 4).Value = 8.760450027789034
$ws.Cells.Item(4, 6).Value = 35.25278025305077
$ws.Cells.Item(4, 7).Value = 3.670441102361348
$ws.Cells.Item(4, 10).Value = 10.42402224957335
$ws.Cells.Item(4, 12).Value = 11.75299255112578
$ws.Cells.Item(4, 14).Value = 18.50069017953968
$ws.Cells.Item(4, 15).Value = 26.78946247861314
$ws.Cells.Item(5, 2).Value = 18.20126434974348
$ws.Cells.Item(5, 3).Value = 8.59071449158276
$ws.Cells.Item(5, 4).Value = 8.761435203341607
$ws.Cells.Item(5, 6).Value = 35.27084891712004
$ws.Cells.Item(5, 7).Value = 3.671056027785638
$ws.Cells.Item(5, 10).Value = 10.43107339011294
$ws.Cells.Item(5, 12).Value = 11.7505623477816
$ws.Cells.Item(5, 14).Value = 18.51588976710829
$ws.Cells.Item(5, 15).Value = 26.80821799217248
$ws.Cells.Item(6, 2).Value = 18.18422682967573
$ws.Cells.Item(6, 3).Value = 8.582079564339207
$ws.Cells.Item(6, 4).Value = 8.761613554097078
$ws.Cells.Item(6, 6).Value = 35.27397916093767
$ws.Cells.Item(6, 7).Value = 3.671159272650876
$ws.Cells.Item(6, 10).Value = 10.43225793617222
$ws.Cells.Item(6, 12).Value = 11.75018242817807
$ws.Cells.Item(6, 14).Value = 18.51844221142537
$ws.Cells.Item(6, 15).Value = 26.81142621792832
$ws.Cells.Item(7, 2).Value = 18.30242460147214
$ws.Cells.Item(7, 3).Value = 8.641747401377188
$ws.Cells.Item(7, 4).Value = 8.760462324967564
$ws.Cells.Item(7, 6).Value = 35.25301521826982
$ws.Cells.Item(7, 7).Value = 3.670449319279049
$ws.Cells.Item(7, 10).Value = 10.42411642500871
$ws.Cells.Item(7, 12).Value = 11.75295819493511
$ws.Cells.Item(7, 14).Value = 18.5008932522995
$ws.Cells.Item(7, 15).Value = 26.78970912538798
$ws.Cells.Item(8, 2).Value = 18.82148302427738
$ws.Cells.Item(8, 3).Value = 8.898089863880886
$ws.Cells.Item(8, 4).Value = 8.757522943053694
$ws.Cells.Item(8, 6).Value = 35.17942688400883
$ws.Cells.Item(8, 7).Value = 3.667481298695914
$ws.Cells.Item(8, 10).Value = 10.39017865241394
$ws.Cells.Item(8, 12).Value = 11.76862345148812
$ws.Cells.Item(8, 14).Value = 18.42759792126768
$ws.Cells.Item(8, 15).Value = 26.70755081847898
$ws.Cells.Item(9, 2).Value = 19.82111640085552
$ws.Cells.Item(9, 3).Value = 9.372875578738551
$ws.Cells.Item(9, 4).Value = 8.75940108602085
$ws.Cells.Item(9, 6).Value = 35.10336803090799
$ws.Cells.Item(9, 7).Value = 3.66224946738797
$ws.Cells.Item(9, 10).Value = 10.33073822092527
$ws.Cells.Item(9, 12).Value = 11.81145751558219
$ws.Cells.Item(9, 14).Value = 18.29868822700542
$ws.Cells.Item(9, 15).Value = 26.59582593903145
$ws.Cells.Item(10, 2).Value = 20.53675997140988
$ws.Cells.Item(10, 3).Value = 9.702362608169377
$ws.Cells.Item(10, 4).Value = 8.765390916028229
$ws.Cells.Item(10, 6).Value = 35.08925350994475
$ws.Cells.Item(10, 7).Value = 3.658760550079538
$ws.Cells.Item(10, 10).Value = 10.29136754045221
$ws.Cells.Item(10, 12).Value = 11.85015311213044
$ws.Cells.Item(10, 14).Value = 18.21294523217205
$ws.Cells.Item(10, 15).Value = 26.54404156589909
$ws.Cells.Item(11, 2).Value = 20.85675795872734
$ws.Cells.Item(11, 3).Value = 9.847618229514728
$ws.Cells.Item(11, 4).Value = 8.769105335777363
$ws.Cells.Item(11, 6).Value = 35.09191751299467
$ws.Cells.Item(11, 7).Value = 3.657249603091411
$ws.Cells.Item(11, 10).Value = 10.27438312165175
$ws.Cells.Item(11, 12).Value = 11.86928770243605
$ws.Cells.Item(11, 14).Value = 18.17587121259213
$ws.Cells.Item(11, 15).Value = 26.52709558810103
$ws.Cells.Item(12, 2).Value = 20.97702062242284
$ws.Cells.Item(12, 3).Value = 9.901924075588241
$ws.Cells.Item(12, 4).Value = 8.770653104124037
$ws.Cells.Item(12, 6).Value = 35.09423230810105
$ws.Cells.Item(12, 7).Value = 3.656688339472589
$ws.Cells.Item(12, 10).Value = 10.26808409669232
$ws.Cells.Item(12, 12).Value = 11.87675017927269
$ws.Cells.Item(12, 14).Value = 18.16210886505385
$ws.Cells.Item(12, 15).Value = 26.5216311449524
$ws.Cells.Item(13, 2).Value = 20.95116229207878
$ws.Cells.Item(13, 3).Value = 9.890259947565546
$ws.Cells.Item(13, 4).Value = 8.770313502958293
$ws.Cells.Item(13, 6).Value = 35.09367570765855
$ws.Cells.Item(13, 7).Value = 3.65680873376488
$ws.Cells.Item(13, 10).Value = 10.26943481432201
$ws.Cells.Item(13, 12).Value = 11.87513342701536
$ws.Cells.Item(13, 14).Value = 18.16506053680554
$ws.Cells.Item(13, 15).Value = 26.52276561559797
$ws.Cells.Item(14, 2).Value = 20.86667098338548
$ws.Cells.Item(14, 3).Value = 9.852100201944751
$ws.Cells.Item(14, 4).Value = 8.769229848342901
$ws.Cells.Item(14, 6).Value = 35.09208178419404
$ws.Cells.Item(14, 7).Value = 3.657203209474149
$ws.Cells.Item(14, 10).Value = 10.27386224218725
$ws.Cells.Item(14, 12).Value = 11.86989732498879
$ws.Cells.Item(14, 14).Value = 18.17473343388196
$ws.Cells.Item(14, 15).Value = 26.52662692329265
$ws.Cells.Item(15, 2).Value = 20.81479534806174
$ws.Cells.Item(15, 3).Value = 9.828634235644044
$ws.Cells.Item(15, 4).Value = 8.768584432894436
$ws.Cells.Item(15, 6).Value = 35.0912755109684
$ws.Cells.Item(15, 7).Value = 3.657446255003619
$ws.Cells.Item(15, 10).Value = 10.27659142464247
$ws.Cells.Item(15, 12).Value = 11.86671816065085
$ws.Cells.Item(15, 14).Value = 18.18069438166819
$ws.Cells.Item(15, 15).Value = 26.5291161970144
$ws.Cells.Item(16, 2).Value = 20.51572559998337
$ws.Cells.Item(16, 3).Value = 9.692773719894864
$ws.Cells.Item(16, 4).Value = 8.765167990267136
$ws.Cells.Item(16, 6).Value = 35.08926220946888
$ws.Cells.Item(16, 7).Value = 3.658860822203537
$ws.Cells.Item(16, 10).Value = 10.29249609443341
$ws.Cells.Item(16, 12).Value = 11.8489331375817
$ws.Cells.Item(16, 14).Value = 18.21540688524964
$ws.Cells.Item(16, 15).Value = 26.54528223536958
$ws.Cells.Item(17, 2).Value = 20.33074953513401
$ws.Cells.Item(17, 3).Value = 9.608217180807502
$ws.Cells.Item(17, 4).Value = 8.763324808807758
$ws.Cells.Item(17, 6).Value = 35.09035397966886
$ws.Cells.Item(17, 7).Value = 3.659748085193924
$ws.Cells.Item(17, 10).Value = 10.30248979371058
$ws.Cells.Item(17, 12).Value = 11.83841239223965
$ws.Cells.Item(17, 14).Value = 18.23719579657669
$ws.Cells.Item(17, 15).Value = 26.556894367252
$ws.Cells.Item(18, 2).Value = 20.22384002988365
$ws.Cells.Item(18, 3).Value = 9.559148560977874
$ws.Cells.Item(18, 4).Value = 8.762357876981403
$ws.Cells.Item(18, 6).Value = 35.09183712324398
$ws.Cells.Item(18, 7).Value = 3.660265589000166
$ws.Cells.Item(18, 10).Value = 10.30832504607913
$ws.Cells.Item(18, 12).Value = 11.83250554529518
$ws.Cells.Item(18, 14).Value = 18.24990999862441
$ws.Cells.Item(18, 15).Value = 26.56419552750541
$ws.Cells.Item(19, 2).Value = 20.18755729695957
$ws.Cells.Item(19, 3).Value = 9.542461307687143
$ws.Cells.Item(19, 4).Value = 8.762046533415226
$ws.Cells.Item(19, 6).Value = 35.09248616491143
$ws.Cells.Item(19, 7).Value = 3.66044204066085
$ws.Cells.Item(19, 10).Value = 10.3103157435278
$ws.Cells.Item(19, 12).Value = 11.83053049847964
$ws.Cells.Item(19, 14).Value = 18.25424605997522
$ws.Cells.Item(19, 15).Value = 26.56677436510511
$ws.Cells.Item(20, 2).Value = 20.35049483401956
$ws.Cells.Item(20, 3).Value = 9.617263531076343
$ws.Cells.Item(20, 4).Value = 8.763511378809429
$ws.Cells.Item(20, 6).Value = 35.09014924993356
$ws.Cells.Item(20, 7).Value = 3.65965289249359
$ws.Cells.Item(20, 10).Value = 10.30141693210486
$ws.Cells.Item(20, 12).Value = 11.83951742370696
$ws.Cells.Item(20, 14).Value = 18.23485752095834
$ws.Cells.Item(20, 15).Value = 26.55559382796137
$ws.Cells.Item(21, 2).Value = 20.89151380433818
$ws.Cells.Item(21, 3).Value = 9.863327877308942
$ws.Cells.Item(21, 4).Value = 8.769544321013688
$ws.Cells.Item(21, 6).Value = 35.0925145217595
$ws.Cells.Item(21, 7).Value = 3.65708704703245
$ws.Cells.Item(21, 10).Value = 10.27255820421257
$ws.Cells.Item(21, 12).Value = 11.871429445965
$ws.Cells.Item(21, 14).Value = 18.17188476474394
$ws.Cells.Item(21, 15).Value = 26.52546689608439
$ws.Cells.Item(22, 2).Value = 21.2397277156709
$ws.Cells.Item(22, 3).Value = 10.0200553803368
$ws.Cells.Item(22, 4).Value = 8.774309652990304
$ws.Cells.Item(22, 6).Value = 35.10167182671345
$ws.Cells.Item(22, 7).Value = 3.655473622330938
$ws.Cells.Item(22, 10).Value = 10.25447006707285
$ws.Cells.Item(22, 12).Value = 11.89354661969574
$ws.Cells.Item(22, 14).Value = 18.13234131571138
$ws.Cells.Item(22, 15).Value = 26.51133044596521
$ws.Cells.Item(23, 2).Value = 21.05440710834889
$ws.Cells.Item(23, 3).Value = 9.936791572697052
$ws.Cells.Item(23, 4).Value = 8.771691432097349
$ws.Cells.Item(23, 6).Value = 35.09608833293733
$ws.Cells.Item(23, 7).Value = 3.656328945009201
$ws.Cells.Item(23, 10).Value = 10.26405350033238
$ws.Cells.Item(23, 12).Value = 11.88162815424805
$ws.Cells.Item(23, 14).Value = 18.153299104708
$ws.Cells.Item(23, 15).Value = 26.51836668333208
$ws.Cells.Item(24, 2).Value = 20.34156973318771
$ws.Cells.Item(24, 3).Value = 9.61317509373227
$ws.Cells.Item(24, 4).Value = 8.763426741572982
$ws.Cells.Item(24, 6).Value = 35.09023914346094
$ws.Cells.Item(24, 7).Value = 3.659695906044353
$ws.Cells.Item(24, 10).Value = 10.30190169326073
$ws.Cells.Item(24, 12).Value = 11.83901739742454
$ws.Cells.Item(24, 14).Value = 18.23591407133231
$ws.Cells.Item(24, 15).Value = 26.55617985446299
$ws.Cells.Item(25, 2).Value = 19.55345904330934
$ws.Cells.Item(25, 3).Value = 9.247671643820787
$ws.Cells.Item(25, 4).Value = 8.758079668304982
$ws.Cells.Item(25, 6).Value = 35.11661748493612
$ws.Cells.Item(25, 7).Value = 3.663602215912405
$ws.Cells.Item(25, 10).Value = 10.34606075687793
$ws.Cells.Item(25, 12).Value = 11.81145751558219
$ws.Cells.Item(25, 14).Value = 18.29868822700542
$ws.Cells.Item(25, 15).Value = 26.52276561559797

Write-Output "Done updating 216 cells"